$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 45189
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 100112003
$ws.Cells.Item(5, 7).Value = "Ajo"
$ws.Cells.Item(5, 8).Value = "Chino"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 750
$ws.Cells.Item(5, 11).Value = 23000
$ws.Cells.Item(5, 12).Value = 24000
$ws.Cells.Item(5, 13).Value = 23467
$ws.Cells.Item(5, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(5, 15).Value = "China"
$ws.Cells.Item(5, 16).Value = 2347
$ws.Cells.Item(5, 17).Value = 10
$ws.Cells.Item(5, 18).Value = "Hortaliza"
